$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.928.44"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "2.493.09"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'535.04"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").Value = "'136.51"
$ws.Range("E6").Value = "  -2.50%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.567"
$ws.Range("E8").Value = "  +0.42%  "
$ws.Range("D9").Value = "2.514.08"
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("E10").Value = "  -0.94%  "
$ws.Range("E11").Value = "  -2.43%  "
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("D13").Value = "'0.346"
$ws.Range("E13").Value = "  -3.61%  "
$ws.Range("D14").Value = "2.938.06"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").Value = "58.771.73"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("D16").Value = "'22.85"
$ws.Range("E16").Value = "  -2.81%  "
$ws.Range("E17").Value = "  -2.20%  "
$ws.Range("D18").Value = "2.509.70"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").Value = "'11.05"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("D21").Value = "'322.66"
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("D24").Value = "'65.50"
$ws.Range("E24").Value = "  +3.49%  "
$ws.Range("E25").Value = "  -1.33%  "
$ws.Range("E26").Value = "  -1.87%  "
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("D28").Value = "'7.53"
$ws.Range("E28").Value = "  -3.96%  "
$ws.Range("E29").Value = "  -3.88%  "
$ws.Range("D30").Value = "0.0₃0765"
$ws.Range("E30").Value = "  -2.48%  "
$ws.Range("E31").Value = "  -1.82%  "
$ws.Range("D32").Value = "'166.41"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("D33").Value = "'1.17"
$ws.Range("E33").Value = "  +3.97%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("D36").Value = "'18.42"
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("D37").Value = "'4.08"
$ws.Range("E37").Value = "  -4.63%  "
$ws.Range("E38").Value = "  -3.93%  "
$ws.Range("D39").Value = "'36.65"
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("E40").Value = "  -0.36%  "
$ws.Range("D41").Value = "'3.60"
$ws.Range("E41").Value = "  -2.81%  "
$ws.Range("D42").Value = "'283.19"
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("D43").Value = "'5.12"
$ws.Range("E43").Value = "  -2.73%  "
$ws.Range("D44").Value = "'0.996"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "'131.12"
$ws.Range("E45").Value = "  +5.61%  "
$ws.Range("D46").Value = "'0.604"
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D47").Value = "'10.87"
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("E48").Value = "  -1.36%  "
$ws.Range("E49").Value = "  -2.39%  "
$ws.Range("E50").Value = "  -2.89%  "
$ws.Range("D51").Value = "'17.17"
$ws.Range("E51").Value = "  -4.25%  "
